# Apply the "Updated symbol list" data refresh (Wed Dec 14 22:12:01 UTC 2022, GitHub Actions):
#  - Column D (Price) is refreshed with newly scraped quotes for the affected coins.
#  - Column G (Hora) moves from hour "21" to hour "22" for every data row (2-51).
# Cells in this sheet are stored as text (see original t="inlineStr" cells), so we force
# the Text number format before writing each value to avoid Excel silently re-typing the
# value as a number (which would also introduce floating point rounding noise).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "267.58"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "22"

# Row 3
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "22"

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "6.321"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "22"

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06197"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "22"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.598"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "22"

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.683"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "22"

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.389"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "22"

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8340"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "22"

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.01364"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "22"

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1603"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "22"

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08258"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "22"

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03404"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "22"

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03151"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "22"

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09288"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "22"

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.922"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "22"

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001709"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "22"

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.04842"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "22"

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006265"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "22"

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.005372"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "22"

# Row 21
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "22"

# Row 22
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "22"

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.775"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "22"

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.369"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "22"

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3349"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "22"

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1213"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "22"

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0002684"
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "22"

# Row 28
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "22"

# Row 29
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "22"

# Row 30
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "22"

# Row 31
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "22"

# Row 32
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "22"

# Row 33
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "22"

# Row 34
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "22"

# Row 35
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "22"

# Row 36
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "22"

# Row 37
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "22"

# Row 38
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "22"

# Row 39
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "22"

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04655"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "22"

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006929"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "22"

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1154"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "22"

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003351"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "22"

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01229"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "22"

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006233"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "22"

# Row 46
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "22"

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.7005"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "22"

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1684"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "22"

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002101"
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "22"

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.01241"
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "22"

# Row 51
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "22"
